$wb = $excel.ActiveWorkbook

# Update OFF sheet (row 2: Home) - Short Att, Short Comp, Deep Att, Deep Comp
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 344
$wsOff.Range("C2").Value = 249
$wsOff.Range("D2").Value = 174
$wsOff.Range("E2").Value = 81

# Update DEF sheet (row 2: Home) - Short Att, Short Comp, Deep Att, Deep Comp, Short Int
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 386
$wsDef.Range("C2").Value = 264
$wsDef.Range("D2").Value = 91
$wsDef.Range("E2").Value = 45
$wsDef.Range("F2").Value = 13

$wb.Save()
